$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column (D) retains text formatting, since several values
# look like ambiguous numbers (e.g. "1.004") and would otherwise be
# auto-converted to numeric by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "28.966.04"
$ws.Range("E2").Value = "  +1.89%  "

$ws.Range("D3").Value = "1.896.31"
$ws.Range("E3").Value = "  +1.51%  "

$ws.Range("D4").Value = "1.004"
$ws.Range("E4").Value = "  +0.22%  "

$ws.Range("D5").Value = "332.87"
$ws.Range("E5").Value = "  -1.43%  "

$ws.Range("D6").Value = "1.003"
$ws.Range("E6").Value = "  +0.14%  "

$ws.Range("D7").Value = "0.4636"
$ws.Range("E7").Value = "  -1.10%  "

$ws.Range("D8").Value = "0.4107"
$ws.Range("E8").Value = "  +3.34%  "

$ws.Range("D9").Value = "47.46"
$ws.Range("E9").Value = "  -0.61%  "

$ws.Range("D10").Value = "0.07968"
$ws.Range("E10").Value = "  -0.73%  "

$ws.Range("D11").Value = "1.002"
$ws.Range("E11").Value = "  +0.30%  "

$ws.Range("D12").Value = "21.77"
$ws.Range("E12").Value = "  -0.92%  "

$ws.Range("D13").Value = "1.899.32"
$ws.Range("E13").Value = "  +1.82%  "

$ws.Range("D14").Value = "5.919"
$ws.Range("E14").Value = "  -2.13%  "

$ws.Range("D15").Value = "7.070"
$ws.Range("E15").Value = "  -2.66%  "

$ws.Range("B16").Value = "BinanceUSD"
$ws.Range("C16").Value = "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
$ws.Range("D16").Value = "1.002"
$ws.Range("E16").Value = "  -0.04%  "

$ws.Range("B17").Value = "Litecoin"
$ws.Range("C17").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D17").Value = "89.02"
$ws.Range("E17").Value = "  -1.80%  "

$ws.Range("D18").Value = "0.00001031"
$ws.Range("E18").Value = "  -0.81%  "

$ws.Range("D19").Value = "0.06565"
$ws.Range("E19").Value = "  -0.78%  "

$ws.Range("D20").Value = "17.49"
$ws.Range("E20").Value = "  -0.24%  "

$ws.Range("D21").Value = "1.001"
$ws.Range("E21").Value = "  -0.01%  "

$ws.Range("D22").Value = "29.036.62"
$ws.Range("E22").Value = "  +2.08%  "

$ws.Range("D23").Value = "5.427"
$ws.Range("E23").Value = "  -0.89%  "

$ws.Range("D24").Value = "11.26"
$ws.Range("E24").Value = "  +1.90%  "

$ws.Range("D25").Value = "2.212"
$ws.Range("E25").Value = "  -2.59%  "

$ws.Range("D26").Value = "2.125.61"
$ws.Range("E26").Value = "  +1.84%  "

$ws.Range("D27").Value = "157.21"
$ws.Range("E27").Value = "  -2.21%  "

$ws.Range("D28").Value = "19.63"
$ws.Range("E28").Value = "  -0.68%  "

$ws.Range("D29").Value = "2.122"
$ws.Range("E29").Value = "  +0.20%  "

$ws.Range("D30").Value = "5.419"
$ws.Range("E30").Value = "  -1.25%  "

$ws.Range("D31").Value = "117.87"
$ws.Range("E31").Value = "  -1.83%  "

$ws.Range("D32").Value = "0.9795"
$ws.Range("E32").Value = "  +0.94%  "

$ws.Range("D33").Value = "0.09397"
$ws.Range("E33").Value = "  -1.22%  "

$ws.Range("D34").Value = "1.425"
$ws.Range("E34").Value = "  +3.66%  "

$ws.Range("D35").Value = "3.598"
$ws.Range("E35").Value = "  +0.32%  "

$ws.Range("D36").Value = "5.280"
$ws.Range("E36").Value = "  -1.44%  "

$ws.Range("D37").Value = "0.06075"
$ws.Range("E37").Value = "  -0.35%  "

$ws.Range("D38").Value = "0.02239"
$ws.Range("E38").Value = "  -0.42%  "

$ws.Range("D39").Value = "8.359"
$ws.Range("E39").Value = "  +0.47%  "

$ws.Range("D40").Value = "1.170"
$ws.Range("E40").Value = "  -0.90%  "

$ws.Range("D41").Value = "1.001"
$ws.Range("E41").Value = "  +0.01%  "

$ws.Range("D42").Value = "0.5788"
$ws.Range("E42").Value = "  -2.44%  "

$ws.Range("D43").Value = "10.17"
$ws.Range("E43").Value = "  -1.33%  "

$ws.Range("D44").Value = "0.1817"
$ws.Range("E44").Value = "  -3.10%  "

$ws.Range("D45").Value = "1.262"
$ws.Range("E45").Value = "  -1.41%  "

$ws.Range("D46").Value = "2.295"
$ws.Range("E46").Value = "  +11.12%  "

$ws.Range("B47").Value = "Decentraland"
$ws.Range("C47").Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Range("D47").Value = "0.5500"
$ws.Range("E47").Value = "  -1.03%  "

$ws.Range("B48").Value = "EnergySwap"
$ws.Range("C48").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D48").Value = "12.04"
$ws.Range("E48").Value = "  -0.97%  "

$ws.Range("D49").Value = "1.907"
$ws.Range("E49").Value = "  -2.45%  "

$ws.Range("D50").Value = "0.07032"
$ws.Range("E50").Value = "  -3.42%  "

$ws.Range("D51").Value = "46.77"
$ws.Range("E51").Value = "  +18.18%  "
